# Fixes across the board
# - Removes the erroneous "quest item" data (quest_item_name / quest_item_drop_chance,
#   columns AM/AN) that had been mistakenly filled in on monsters that should not
#   drop a quest item (row 6, and rows 10-15), leaving it populated only for the
#   rows that actually do drop a quest item (rows 7-9).
# - Removes the extra blank trailing row (row 16) that had been added under the
#   data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the erroneous quest_item_name / quest_item_drop_chance values.
$ws.Range("AM6:AN6").ClearContents()
$ws.Range("AM10:AN15").ClearContents()

# Remove the now-unnecessary blank row at the bottom of the sheet.
$ws.Rows.Item(16).Delete()

# Reflect the resulting view/selection position.
[void]$ws.Range("AM10:AN15").Select()
